$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-10 22:18:52"

$ws.Range("E3").Value = "2026-02-10 22:18:54"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "96%"
$ws.Range("G3").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("I3").Value = "25.8 mm"
$ws.Range("L3").Value = "55.4 km/h - 254º 21:55 TU"

$ws.Range("E4").Value = "2026-02-10 22:18:57"

$ws.Range("E5").Value = "2026-02-10 22:18:59"

$ws.Range("E6").Value = "2026-02-10 22:19:02"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "94%"
$ws.Range("G6").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("J6").Value = "1004.0 hPa"

$ws.Range("E7").Value = "2026-02-10 22:19:05"

$ws.Range("E8").Value = "2026-02-10 22:19:07"
$ws.Range("O8").Value = "12.4 °C"

$ws.Range("E9").Value = "2026-02-10 22:19:10"

$ws.Range("E10").Value = "2026-02-10 22:19:13"

$ws.Range("E11").Value = "2026-02-10 22:19:15"

$ws.Range("E12").Value = "2026-02-10 22:19:18"
$ws.Range("I12").Value = "4.4 mm"

$ws.Range("E13").Value = "2026-02-10 22:19:20"

$ws.Range("E14").Value = "2026-02-10 22:19:23"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "88%"
$ws.Range("G14").Copy()
$ws.Range("H14").PasteSpecial(-4122)

$ws.Range("E15").Value = "2026-02-10 22:19:25"
$ws.Range("O15").Value = "9.2 °C"

$ws.Range("E16").Value = "2026-02-10 22:19:28"
$ws.Range("I16").Value = "26.7 mm"

$ws.Range("E17").Value = "2026-02-10 22:19:31"
$ws.Range("O17").Value = "4.9 °C"

$ws.Range("E18").Value = "2026-02-10 22:19:33"
$ws.Range("J18").Value = "1004.1 hPa"

$ws.Range("E19").Value = "2026-02-10 22:19:36"

$ws.Range("E20").Value = "2026-02-10 22:19:38"
$ws.Range("I20").Value = "12.0 mm"

$ws.Range("E21").Value = "2026-02-10 22:19:41"
$ws.Range("J21").Value = "1005.9 hPa"
$ws.Range("O21").Value = "7.4 °C"

$ws.Range("E22").Value = "2026-02-10 22:19:44"

$ws.Range("E23").Value = "2026-02-10 22:19:46"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "91%"
$ws.Range("G23").Copy()
$ws.Range("H23").PasteSpecial(-4122)

$ws.Range("E24").Value = "2026-02-10 22:19:49"

$ws.Range("E25").Value = "2026-02-10 22:19:52"
$ws.Range("L25").Value = "49.7 km/h - 241º 21:57 TU"

$ws.Range("E26").Value = "2026-02-10 22:19:54"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "84%"
$ws.Range("G26").Copy()
$ws.Range("H26").PasteSpecial(-4122)

$ws.Range("E27").Value = "2026-02-10 22:19:57"

$ws.Range("E28").Value = "2026-02-10 22:19:59"

$ws.Range("E29").Value = "2026-02-10 22:20:02"

$ws.Range("E30").Value = "2026-02-10 22:20:04"

$ws.Range("E31").Value = "2026-02-10 22:20:07"
$ws.Range("J31").Value = "1003.3 hPa"
$ws.Range("O31").Value = "10.6 °C"

$ws.Range("E32").Value = "2026-02-10 22:20:09"

$ws.Range("E33").Value = "2026-02-10 22:20:12"
$ws.Range("J33").Value = "1006.2 hPa"

$ws.Range("E34").Value = "2026-02-10 22:20:14"

$ws.Range("E35").Value = "2026-02-10 22:20:17"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = "68%"
$ws.Range("G35").Copy()
$ws.Range("H35").PasteSpecial(-4122)

$ws.Range("E36").Value = "2026-02-10 22:20:19"
$ws.Range("J36").Value = "1004.2 hPa"

$ws.Range("E37").Value = "2026-02-10 22:20:22"

$ws.Range("E38").Value = "2026-02-10 22:20:25"

$ws.Range("E39").Value = "2026-02-10 22:20:27"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "80%"
$ws.Range("G39").Copy()
$ws.Range("H39").PasteSpecial(-4122)
$ws.Range("I39").Value = "11.8 mm"
$ws.Range("L39").Value = "66.6 km/h - 297º 21:47 TU"

$ws.Range("E40").Value = "2026-02-10 22:20:30"

$ws.Range("E41").Value = "2026-02-10 22:20:32"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "79%"
$ws.Range("G41").Copy()
$ws.Range("H41").PasteSpecial(-4122)
$ws.Range("L41").Value = "44.6 km/h - 298º 21:55 TU"
$ws.Range("O41").Value = "14.9 °C"

$ws.Range("E42").Value = "2026-02-10 22:20:35"

$ws.Range("E43").Value = "2026-02-10 22:20:37"
$ws.Range("O43").Value = "10.1 °C"

$ws.Range("E44").Value = "2026-02-10 22:20:39"
$ws.Range("L44").Value = "56.9 km/h - 185º 21:44 TU"

$ws.Range("E45").Value = "2026-02-10 22:20:42"

$ws.Range("E46").Value = "2026-02-10 22:20:44"
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "78%"
$ws.Range("G46").Copy()
$ws.Range("H46").PasteSpecial(-4122)
$ws.Range("O46").Value = "15.1 °C"

$excel.CutCopyMode = $false